# Fidi watchlist refresh: rewrite rows 2..39 (cols A-F) of Sheet1 with the
# new NSE ticker lists (Buying Opportunity / support Zone / long buildup /
# Short buildup / FII ENTERING), extending the used range from A1:F18 to
# A1:F39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 39
$firstDataRow = 2
$oldLastRow = 18

# --- extend column A (the bordered/bold "#" index column) down to row 39,
#     carrying forward the same style used by A2:A18 -------------------------
if ($lastRow -gt $oldLastRow) {
    $ws.Range("A$firstDataRow").Copy() | Out-Null
    $ws.Range("A" + ($oldLastRow + 1) + ":A$lastRow").PasteSpecial(-4122) | Out-Null
}

$rowCount = $lastRow - $firstDataRow + 1
$aValues = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $aValues[$i,0] = $i
}
$ws.Range("A$firstDataRow" + ":A$lastRow").Value = $aValues

# --- rewrite columns B:F for every data row ---------------------------------
$data = New-Object 'object[,]' $rowCount,5

$data[0,0] = "NSE:ABMINTLLTD"
$data[0,1] = "NSE:ATALREAL"
$data[0,2] = "NSE:BHARATFORG"
$data[0,4] = "NSE:ACC"

$data[1,0] = "NSE:ACC"
$data[1,1] = "NSE:BBOX"
$data[1,2] = "NSE:COALINDIA"
$data[1,4] = "NSE:BHARATFORG"

$data[2,0] = "NSE:ASPINWALL"
$data[2,1] = "NSE:BFINVEST"
$data[2,2] = "NSE:DIVISLAB"
$data[2,4] = "NSE:BHEL"

$data[3,0] = "NSE:BALMLAWRIE"
$data[3,1] = "NSE:DVL"
$data[3,2] = "NSE:HEROMOTOCO"
$data[3,4] = "NSE:DIVISLAB"

$data[4,0] = "NSE:BASML"
$data[4,1] = "NSE:EPIGRAL"
$data[4,2] = "NSE:ICICIBANK"
$data[4,4] = "NSE:HEROMOTOCO"

$data[5,0] = "NSE:BHARATFORG"
$data[5,1] = "NSE:FOCUS"
$data[5,2] = "NSE:INDUSINDBK"
$data[5,4] = "NSE:RECLTD"

$data[6,0] = "NSE:BHEL"
$data[6,1] = "NSE:GROBTEA"
$data[6,2] = "NSE:IRCTC"

$data[7,0] = "NSE:BIKAJI"
$data[7,1] = "NSE:GSPL"
$data[7,2] = "NSE:JINDALSTEL"

$data[8,0] = "NSE:COMPUSOFT"
$data[8,1] = "NSE:HIKAL"
$data[8,2] = "NSE:PETRONET"

$data[9,0] = "NSE:CONCOR"
$data[9,1] = "NSE:JUSTDIAL"

$data[10,0] = "NSE:CREATIVE"
$data[10,1] = "NSE:LAXMIMACH"

$data[11,0] = "NSE:DALMIASUG"
$data[11,1] = "NSE:M&MFIN"

$data[12,0] = "NSE:DHAMPURSUG"
$data[12,1] = "NSE:MAHEPC"

$data[13,0] = "NSE:DIVISLAB"
$data[13,1] = "NSE:NDGL"

$data[14,0] = "NSE:EICHERMOT"
$data[14,1] = "NSE:PLAZACABLE"

$data[15,0] = "NSE:FINEORG"
$data[15,1] = "NSE:RAMRAT"

$data[16,0] = "NSE:GPTINFRA"
$data[16,1] = "NSE:RPOWER"

$data[17,0] = "NSE:HDFCSENSEX"

$data[18,0] = "NSE:HEROMOTOCO"

$data[19,0] = "NSE:IPL"

$data[20,0] = "NSE:JYOTHYLAB"

$data[21,0] = "NSE:KIRLOSBROS"

$data[22,0] = "NSE:KSL"

$data[23,0] = "NSE:LEMONTREE"

$data[24,0] = "NSE:MAITHANALL"

$data[25,0] = "NSE:MALLCOM"

$data[26,0] = "NSE:MASKINVEST"

$data[27,0] = "NSE:MOVALUE"

$data[28,0] = "NSE:NIF100BEES"

$data[29,0] = "NSE:ORBTEXP"

$data[30,0] = "NSE:PANACEABIO"

$data[31,0] = "NSE:PEL"

$data[32,0] = "NSE:PNC"

$data[33,0] = "NSE:RADIOCITY"

$data[34,0] = "NSE:RECLTD"

$data[35,0] = "NSE:REDINGTON"

$data[36,0] = "NSE:RELAXO"

$data[37,0] = "NSE:ROML"

$ws.Range("B$firstDataRow" + ":F$lastRow").Value = $data
